$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Insert the "Vi förväntar oss..." paragraph right after the
#    "Nedan presenteras..." paragraph near the top of the document.
# ------------------------------------------------------------------
$anchor = $d.Content
$found = $anchor.Find.Execute(
    "Nedan presenteras fynd av naturvårdsarter och fridlysta arter som gjorts i det avverkningsanmälda området, samt relevanta utdrag ur standarderna för FSC, Chain of Custody, Controlled Wood och PEFC. I BILAGA 1 finns artfakta om fridlysta arter.",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

$insertedEnd = -1
if ($found) {
    $anchor.Collapse(0)            # wdCollapseEnd
    $anchor.InsertParagraphAfter()
    $newParaRange = $d.Range($anchor.End + 1, $anchor.End + 1)
    $newParaRange.Text = "Vi förväntar oss att ni återkommer med ett skriftligt svar på vårt klagomål och även beskriver vilka korrigerande åtgärder ni satt in för att rätta till identifierade brister i er efterlevnad av den svenska FSC standarden."
    $insertedEnd = $newParaRange.End
}

# ------------------------------------------------------------------
# 2) Remove the old copy of that paragraph further down (just before
#    the page break to "BILAGA 1"), together with the two empty
#    paragraphs that preceded it. Search strictly AFTER the text we
#    just inserted above so we land on the original occurrence.
# ------------------------------------------------------------------
$searchStart = $insertedEnd
if ($searchStart -lt 0) { $searchStart = 0 }
$old = $d.Range($searchStart, $d.Content.End)
$foundOld = $old.Find.Execute(
    "Vi förväntar oss att ni återkommer med ett skriftligt svar på vårt klagomål och även beskriver vilka korrigerande åtgärder ni satt in för att rätta till identifierade brister i er efterlevnad av den svenska FSC standarden.",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($foundOld) {
    # Determine the 1-based paragraph index of the found (old) paragraph.
    $idx = $d.Range(0, $old.Start).Paragraphs.Count + 1

    $startDeletePara = $d.Paragraphs.Item($idx - 2)   # first of the two empty paragraphs
    $endDeletePara = $d.Paragraphs.Item($idx)         # the old "Vi förväntar oss..." paragraph

    $deleteRange = $d.Range($startDeletePara.Range.Start, $endDeletePara.Range.End)
    $deleteRange.Delete()
}

# ------------------------------------------------------------------
# 3) Update the date in the first-page header from 2023-11-13 to
#    2023-11-14.
# ------------------------------------------------------------------
foreach ($sec in $d.Sections) {
    $hdr = $sec.Headers.Item(2)   # wdHeaderFooterFirstPage
    if ($hdr.Exists) {
        $hdr.Range.Find.Execute("2023-11-13", $true, $false, $false, $false, $false, $true, 1, $false, "2023-11-14", 2)
    }
}
